# Client: Add MoveState Stamina Control System
# "스테미나가 없으면 달리지 못하게 수정"
# Rebalance the movement-speed "Value" column on the Move sheet so that
# walking/running/dashing (and their held-weapon variants) cost more,
# tying movement speed to the new stamina system. Also leaves the Move
# sheet as the active/selected sheet, with D11 selected, matching the
# state the workbook was saved in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Move")

$ws.Range("C3").Value = 0.5    # Ch_Walk
$ws.Range("C4").Value = 0.8    # Ch_Run
$ws.Range("C6").Value = 1      # Ch_Dash
$ws.Range("C7").Value = 0.4    # Ch_HoldWalk
$ws.Range("C8").Value = 0.7    # Ch_HoldRun

# Leave the Move sheet active/selected, matching the saved workbook state.
$ws.Activate()
$ws.Range("D11").Select() | Out-Null
